$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("W7").Value = 2.75
$ws.Range("X7").Value = 1.44

# Row 8
$ws.Range("I8").Value = 3.6
$ws.Range("K8").Value = 2.25
$ws.Range("L8").Value = 4
$ws.Range("O8").Value = 1.22
$ws.Range("P8").Value = 4
$ws.Range("S8").Value = 1.75
$ws.Range("T8").Value = 2.05
$ws.Range("W8").Value = 2.75
$ws.Range("X8").Value = 1.4
$ws.Range("Y8").Value = 1.36
$ws.Range("Z8").Value = 3
$ws.Range("AA8").Value = 1.67
$ws.Range("AB8").Value = 2.1
$ws.Range("AC8").Value = 8.5
$ws.Range("AD8").Value = 10
$ws.Range("AH8").Value = 23
$ws.Range("AI8").Value = 12
$ws.Range("AK8").Value = 13
$ws.Range("AQ8").Value = 29
$ws.Range("AR8").Value = 34
$ws.Range("AS8").Value = 151

# Row 10
$ws.Range("N10").Value = 10
$ws.Range("O10").Value = 1.33
$ws.Range("P10").Value = 3.25
$ws.Range("S10").Value = 2.08
$ws.Range("T10").Value = 1.73
$ws.Range("W10").Value = 3.75
$ws.Range("X10").Value = 1.25

# Row 11
$ws.Range("M11").Value = 1.03
$ws.Range("N11").Value = 15
$ws.Range("O11").Value = 1.2
$ws.Range("P11").Value = 4.33
$ws.Range("S11").Value = 1.67
$ws.Range("T11").Value = 2.15
$ws.Range("W11").Value = 2.63
$ws.Range("X11").Value = 1.44

# Row 12
$ws.Range("G12").Value = 3.6
$ws.Range("H12").Value = 3.4
$ws.Range("I12").Value = 2.05
$ws.Range("K12").Value = 2.1
$ws.Range("L12").Value = 2.75
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 10
$ws.Range("O12").Value = 1.3
$ws.Range("P12").Value = 3.4
$ws.Range("S12").Value = 2.05
$ws.Range("T12").Value = 1.8
$ws.Range("W12").Value = 3.5
$ws.Range("X12").Value = 1.29
$ws.Range("AA12").Value = 1.8
$ws.Range("AB12").Value = 1.95
$ws.Range("AE12").Value = 12
$ws.Range("AI12").Value = 10
$ws.Range("AK12").Value = 15
$ws.Range("AL12").Value = 51
$ws.Range("AM12").Value = 7.5
$ws.Range("AO12").Value = 9
$ws.Range("AP12").Value = 19
$ws.Range("AQ12").Value = 17
$ws.Range("AS12").Value = 251

# Row 17
$ws.Range("G17").Value = 2.52
$ws.Range("H17").Value = 3.6
$ws.Range("I17").Value = 2.45
$ws.Range("J17").Value = 2.9
$ws.Range("K17").Value = 2.37
$ws.Range("L17").Value = 2.92
$ws.Range("P17").Value = 5.1
$ws.Range("S17").Value = 1.42
$ws.Range("T17").Value = 2.67
$ws.Range("AA17").Value = 1.36
$ws.Range("AB17").Value = 2.9
$ws.Range("AC17").Value = 16
$ws.Range("AD17").Value = 19.5
$ws.Range("AG17").Value = 17.5
$ws.Range("AH17").Value = 17.5
$ws.Range("AM17").Value = 14.5
$ws.Range("AN17").Value = 17
$ws.Range("AQ17").Value = 17
$ws.Range("AR17").Value = 18
$ws.Range("AS17").Value = 120

# Row 18
$ws.Range("G18").Value = 3.8
$ws.Range("H18").Value = 3.85
$ws.Range("I18").Value = 1.8
$ws.Range("J18").Value = 4.05
$ws.Range("K18").Value = 2.32
$ws.Range("L18").Value = 2.3
$ws.Range("W18").Value = 2.3
$ws.Range("X18").Value = 1.55
$ws.Range("Z18").Value = 3.25
$ws.Range("AA18").Value = 1.53
$ws.Range("AC18").Value = 14.5
$ws.Range("AD18").Value = 24
$ws.Range("AG18").Value = 30
$ws.Range("AH18").Value = 30
$ws.Range("AJ18").Value = 7.8
$ws.Range("AK18").Value = 12.5
$ws.Range("AM18").Value = 10
$ws.Range("AP18").Value = 16
$ws.Range("AQ18").Value = 12.5
